$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "41.984.34"
$ws.Range("E2").Value = "  -1.18%  "
$ws.Range("D3").Value = "2.256.48"
$ws.Range("E3").Value = "  -1.32%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'306.93"
$ws.Range("E5").Value = "  +0.05%  "
$ws.Range("D6").Value = "'96.54"
$ws.Range("E6").Value = "  -0.16%  "
$ws.Range("D7").Value = "'0.522"
$ws.Range("E7").Value = "  -1.89%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.486"
$ws.Range("E9").Value = "  -1.87%  "
$ws.Range("D10").Value = "'34.91"
$ws.Range("E10").Value = "  -3.86%  "
$ws.Range("D11").Value = "'0.0784"
$ws.Range("E11").Value = "  -2.53%  "
$ws.Range("D12").Value = "'0.113"
$ws.Range("E12").Value = "  +0.29%  "
$ws.Range("D13").Value = "'6.77"
$ws.Range("E13").Value = "  +0.77%  "
$ws.Range("D14").Value = "2.611.69"
$ws.Range("E14").Value = "  -1.16%  "
$ws.Range("D15").Value = "'14.52"
$ws.Range("E15").Value = "  -0.60%  "
$ws.Range("D16").Value = "2.260.90"
$ws.Range("E16").Value = "  -1.46%  "
$ws.Range("D17").Value = "'0.785"
$ws.Range("E17").Value = "  -2.37%  "
$ws.Range("D18").Value = "41.861.65"
$ws.Range("E18").Value = "  -1.23%  "
$ws.Range("D19").Value = "'12.13"
$ws.Range("E19").Value = "  -5.77%  "
$ws.Range("D20").Value = "0.0₃0897"
$ws.Range("E20").Value = "  -2.46%  "
$ws.Range("D21").Value = "'5.93"
$ws.Range("E21").Value = "  -1.31%  "
$ws.Range("D22").Value = "'67.44"
$ws.Range("E22").Value = "  -0.68%  "
$ws.Range("D23").Value = "'235.37"
$ws.Range("E23").Value = "  -3.19%  "
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "'2.56"
$ws.Range("E24").Value = "  -1.48%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "'1.96"
$ws.Range("E25").Value = "  -0.05%  "
$ws.Range("E26").Value = "  -0.09%  "
$ws.Range("D27").Value = "'23.41"
$ws.Range("E27").Value = "  -2.37%  "
$ws.Range("D28").Value = "'36.25"
$ws.Range("E28").Value = "  -0.71%  "
$ws.Range("D29").Value = "'9.46"
$ws.Range("E29").Value = "  -1.39%  "
$ws.Range("D30").Value = "'2.11"
$ws.Range("E30").Value = "  +0.09%  "
$ws.Range("D31").Value = "'163.70"
$ws.Range("E31").Value = "  +1.29%  "
$ws.Range("B32").Value = "FirstDigitalUSD"
$ws.Range("C32").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D32").Value = "'1.00"
$ws.Range("E32").Value = "  +0.09%  "
$ws.Range("B33").Value = "Filecoin"
$ws.Range("C33").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D33").Value = "'5.19"
$ws.Range("E33").Value = "  -2.50%  "
$ws.Range("D34").Value = "'3.12"
$ws.Range("E34").Value = "  +0.17%  "
$ws.Range("D35").Value = "'0.0731"
$ws.Range("E35").Value = "  -2.89%  "
$ws.Range("D36").Value = "'17.31"
$ws.Range("E36").Value = "  -0.25%  "
$ws.Range("D37").Value = "'2.39"
$ws.Range("E37").Value = "  +0.15%  "
$ws.Range("E38").Value = "  -4.78%  "
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D40").Value = "'1.80"
$ws.Range("E40").Value = "  -4.17%  "
$ws.Range("D41").Value = "'4.09"
$ws.Range("E41").Value = "  -2.17%  "
$ws.Range("D42").Value = "'2.28"
$ws.Range("E42").Value = "  -4.19%  "
$ws.Range("D43").Value = "1.944.42"
$ws.Range("E43").Value = "  -2.96%  "
$ws.Range("D44").Value = "'18.63"
$ws.Range("E44").Value = "  -3.43%  "
$ws.Range("D45").Value = "'0.0278"
$ws.Range("E45").Value = "  -2.29%  "
$ws.Range("D46").Value = "'2.90"
$ws.Range("E46").Value = "  -4.35%  "
$ws.Range("D47").Value = "'9.74"
$ws.Range("E47").Value = "  -5.08%  "
$ws.Range("D48").Value = "'53.08"
$ws.Range("E48").Value = "  -1.63%  "
$ws.Range("D49").Value = "2.484.33"
$ws.Range("E49").Value = "  -1.01%  "
$ws.Range("D50").Value = "'91.77"
$ws.Range("E50").Value = "  -0.28%  "
$ws.Range("D51").Value = "'71.08"
$ws.Range("E51").Value = "  -2.15%  "

Write-Host "Applied cryptos update"
